$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr0 = New-Object 'object[,]' 24,4
$arr0[0,0] = 14.79968816316329
$arr0[0,1] = 8.03515003602395
$arr0[0,2] = 6.012036206773072
$arr0[0,3] = 11.78492940526285
$arr0[1,0] = 14.5855400208801
$arr0[1,1] = 8.011172676253452
$arr0[1,2] = 5.897979913608379
$arr0[1,3] = 11.79297818127113
$arr0[2,0] = 14.4558672562244
$arr0[2,1] = 7.996200055437947
$arr0[2,2] = 5.828601884493501
$arr0[2,3] = 11.79999783055597
$arr0[3,0] = 14.40354837634165
$arr0[3,1] = 7.990036644270222
$arr0[3,2] = 5.80053691995563
$arr0[3,3] = 11.80338111513339
$arr0[4,0] = 14.3948944606791
$arr0[4,1] = 7.989009502945478
$arr0[4,2] = 5.795890523921611
$arr0[4,3] = 11.80397448280957
$arr0[5,0] = 14.45515945600063
$arr0[5,1] = 7.996117182987398
$arr0[5,2] = 5.828222494934082
$arr0[5,3] = 11.80004134199421
$arr0[6,0] = 14.72551377733962
$arr0[6,1] = 8.026935230732779
$arr0[6,2] = 5.972599737882843
$arr0[6,3] = 11.78727361294756
$arr0[7,0] = 15.26710074976902
$arr0[7,1] = 8.085338691822425
$arr0[7,2] = 6.259040091902812
$arr0[7,3] = 11.77870407321926
$arr0[8,0] = 15.66799577569386
$arr0[8,1] = 8.126953335988569
$arr0[8,2] = 6.469074813841146
$arr0[8,3] = 11.78241583276241
$arr0[9,0] = 15.85019826301357
$arr0[9,1] = 8.145591842080817
$arr0[9,2] = 6.564045738772656
$arr0[9,3] = 11.78626725605301
$arr0[10,0] = 15.91910295599487
$arr0[10,1] = 8.152606657655696
$arr0[10,2] = 6.599887613907389
$arr0[10,3] = 11.78803550672206
$arr0[11,0] = 15.90426821528996
$arr0[11,1] = 8.151097830709604
$arr0[11,2] = 6.592174428664828
$arr0[11,3] = 11.78764092538361
$arr0[12,0] = 15.85586930484624
$arr0[12,1] = 8.14616982729045
$arr0[12,2] = 6.566997123107264
$arr0[12,3] = 11.78640652993345
$arr0[13,0] = 15.82620970142631
$arr0[13,1] = 8.143145624723621
$arr0[13,2] = 6.551558346800933
$arr0[13,3] = 11.78569073261288
$arr0[14,0] = 15.65607924564564
$arr0[14,1] = 8.125729286496663
$arr0[14,2] = 6.462853303058451
$arr0[14,3] = 11.78220755631216
$arr0[15,0] = 15.55161837728193
$arr0[15,1] = 8.114969259804589
$arr0[15,2] = 6.40826009720494
$arr0[15,3] = 11.78062390534781
$arr0[16,0] = 15.49152347594921
$arr0[16,1] = 8.108752943248163
$arr0[16,2] = 6.376807620786777
$arr0[16,3] = 11.7799167012798
$arr0[17,0] = 15.47117644268829
$arr0[17,1] = 8.106643528885526
$arr0[17,2] = 6.366150654194366
$arr0[17,3] = 11.77971226891838
$arr0[18,0] = 15.56274012889107
$arr0[18,1] = 8.11611753073235
$arr0[18,2] = 6.414077297414112
$arr0[18,3] = 11.78077141689124
$arr0[19,0] = 15.87008824655538
$arr0[19,1] = 8.147618482591074
$arr0[19,2] = 6.574395906912953
$arr0[19,3] = 11.78676070461556
$arr0[20,0] = 16.07039534812819
$arr0[20,1] = 8.167953935174371
$arr0[20,2] = 6.678448285738018
$arr0[20,3] = 11.79248010125596
$arr0[21,0] = 15.9635608795506
$arr0[21,1] = 8.157123986574726
$arr0[21,2] = 6.622992241672941
$arr0[21,3] = 11.78926284314009
$arr0[22,0] = 15.55771210114986
$arr0[22,1] = 8.115598491541999
$arr0[22,2] = 6.411447544172602
$arr0[22,3] = 11.78070409371161
$arr0[23,0] = 15.11977723508764
$arr0[23,1] = 8.069762202302762
$arr0[23,2] = 6.181451368919016
$arr0[23,3] = 11.77926189758929
$ws.Range("B2:E25").Value = $arr0

$arr1 = New-Object 'object[,]' 24,1
$arr1[0,0] = 3.679432670411195
$arr1[1,0] = 3.681582099668389
$arr1[2,0] = 3.682972340633484
$arr1[3,0] = 3.683556655341615
$arr1[4,0] = 3.683654755906459
$arr1[5,0] = 3.682980148839173
$arr1[6,0] = 3.680159199160319
$arr1[7,0] = 3.675183984725948
$arr1[8,0] = 3.671864426559354
$arr1[9,0] = 3.670426404144819
$arr1[10,0] = 3.669892166334323
$arr1[11,0] = 3.670006766290291
$arr1[12,0] = 3.670382245747544
$arr1[13,0] = 3.670613578888208
$arr1[14,0] = 3.671959850286609
$arr1[15,0] = 3.672804163488245
$arr1[16,0] = 3.673296575953523
$arr1[17,0] = 3.673464465427836
$arr1[18,0] = 3.672713582977113
$arr1[19,0] = 3.67027167893141
$arr1[20,0] = 3.668735824487499
$arr1[21,0] = 3.66955005956491
$arr1[22,0] = 3.672754512600295
$arr1[23,0] = 3.676470693830506
$ws.Range("G2:G25").Value = $arr1

$arr2 = New-Object 'object[,]' 24,2
$arr2[0,0] = 11.12006857629008
$arr2[0,1] = 9.813582256100769
$arr2[1,0] = 10.9716487272167
$arr2[1,1] = 9.802035278995129
$arr2[2,0] = 10.88184665043957
$arr2[2,1] = 9.796681164518178
$arr2[3,0] = 10.84562961642719
$arr2[3,1] = 9.794937787757194
$arr2[4,0] = 10.83963991366469
$arr2[4,1] = 9.794674836664356
$arr2[5,0] = 10.88135662643634
$arr2[5,1] = 9.796655874955766
$arr2[6,0] = 11.06864428511614
$arr2[6,1] = 9.809241872827563
$arr2[7,0] = 11.44450199709611
$arr2[7,1] = 9.847597250504329
$arr2[8,0] = 11.72326405202772
$arr2[8,1] = 9.883965704153464
$arr2[9,0] = 11.8501020091017
$arr2[9,1] = 9.902251717003697
$arr2[10,0] = 11.89809167442662
$arr2[10,1] = 9.909423025426465
$arr2[11,0] = 11.88775877503585
$arr2[11,1] = 9.907867641106391
$arr2[12,0] = 11.85405122576414
$arr2[12,1] = 9.902836780799019
$arr2[13,0] = 11.83339770619452
$arr2[13,1] = 9.899787259489687
$arr2[14,0] = 11.71497156844226
$arr2[14,1] = 9.882805411665508
$arr2[15,0] = 11.64229560867942
$arr2[15,1] = 9.872831187359782
$arr2[16,0] = 11.60049955722686
$arr2[16,1] = 9.867258463168248
$arr2[17,0] = 11.58635042536343
$arr2[17,1] = 9.865399941308285
$arr2[18,0] = 11.65003187697099
$arr2[18,1] = 9.873875993756863
$arr2[19,0] = 11.86395341931503
$arr2[19,1] = 9.904307799191384
$arr2[20,0] = 12.00350384303308
$arr2[20,1] = 9.925633344816474
$arr2[21,0] = 11.92906148941495
$arr2[21,1] = 9.914121334370083
$arr2[22,0] = 11.64653434920424
$arr2[22,1] = 9.873403133176806
$arr2[23,0] = 11.34217440996818
$arr2[23,1] = 9.835771291757267
$ws.Range("K2:L25").Value = $arr2

$arr3 = New-Object 'object[,]' 24,2
$arr3[0,0] = 20.68709000859546
$arr3[0,1] = 28.22637119563035
$arr3[1,0] = 20.75019747451624
$arr3[1,1] = 28.27621778182953
$arr3[2,0] = 20.79078225540027
$arr3[2,1] = 28.31240374237086
$arr3[3,0] = 20.8077839420037
$arr3[3,1] = 28.32855069502705
$arr3[4,0] = 20.81063506341026
$arr3[4,1] = 28.33131642256817
$arr3[5,0] = 20.79100966954371
$arr3[5,1] = 28.31261583692753
$arr3[6,0] = 20.70846904177367
$arr3[6,1] = 28.24239891136136
$arr3[7,0] = 20.56112135572824
$arr3[7,1] = 28.14906862827453
$arr3[8,0] = 20.46163072798656
$arr3[8,1] = 28.10765361139803
$arr3[9,0] = 20.41825568105113
$arr3[9,1] = 28.09472593012516
$arr3[10,0] = 20.40210027541771
$arr3[10,1] = 28.0906814377504
$arr3[11,0] = 20.4055676497009
$arr3[11,1] = 28.09151463749095
$arr3[12,0] = 20.41692116547915
$arr3[12,1] = 28.09437613012272
$arr3[13,0] = 20.42391062316367
$arr3[13,1] = 28.09623971094146
$arr3[14,0] = 20.46450320094906
$arr3[14,1] = 28.10861750908372
$arr3[15,0] = 20.48988711919369
$arr3[15,1] = 28.11772581453311
$arr3[16,0] = 20.50466464860045
$arr3[16,1] = 28.12352113409922
$arr3[17,0] = 20.5096985634251
$arr3[17,1] = 28.12557887327183
$arr3[18,0] = 20.48716660849247
$arr3[18,1] = 28.11669862122089
$arr3[19,0] = 20.41357904950848
$arr3[19,1] = 28.0935125418273
$arr3[20,0] = 20.36705733011362
$arr3[20,1] = 28.08331911859144
$arr3[21,0] = 20.39174336746699
$arr3[21,1] = 28.08830554117468
$arr3[22,0] = 20.48839597819034
$arr3[22,1] = 28.11716127460713
$arr3[23,0] = 20.59943754171936
$arr3[23,1] = 28.16955483253426
$ws.Range("N2:O25").Value = $arr3

Write-Host "done"